$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: fixed width instead of autofit-bestfit ---
$ws.Columns.Item(1).ColumnWidth = 17.14

# --- Row 2: taller custom row height (room for the new discount/amount row) ---
$ws.Rows.Item(2).RowHeight = 37.5

# --- A2: keep centered horizontally, also center vertically now the row is taller ---
$ws.Range("A2").VerticalAlignment = -4108

# --- Unmerge the old "TONG CONG" (grand total) label cells ---
$ws.Range("B3:D3").UnMerge()

# Move the "TONG CONG" label into A3 (first column) instead of the merged B3:D3 block
$ws.Range("A3").Value = "TỔNG CỘNG"
$ws.Range("B3").Value = ""

# Give B3, C3, D3 a uniform full box border (like the rest of the table) now that
# they are independent cells instead of a merged block with a stitched-together border
$ws.Range("E3").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("B3").NumberFormat = "@"

$ws.Range("E3").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("C3").NumberFormat = "@"

$ws.Range("E3").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").HorizontalAlignment = -4108

$excel.CutCopyMode = $false

# --- E3: grand total now dynamically sums every row above it (supports inserted rows) ---
$ws.Range("E3").Formula = '=SUM(E2:INDIRECT("E"&ROW()-1))'

# --- reflect the last place the user clicked before saving ---
$ws.Range("D4").Select()
